$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column in H1, copying formatting from the
# neighboring header cell (G1) so it reuses the existing header style.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Add the data value for the new column in H2
$ws.Range("H2").Value = 1
